$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 15000
$ws.Range("I10").Value = 15000
$ws.Range("K10").Value = 15000
$ws.Range("M10").Value = -14707

$ws.Range("H38").Value = 92.875
$ws.Range("I38").Value = 92.875
$ws.Range("K38").Value = 278.625
$ws.Range("M38").Value = 93.375

$ws.Range("H40").Value = 4521.0444
$ws.Range("I40").Value = 2676.5
$ws.Range("J40").Value = 4804.8203
$ws.Range("K40").Value = 2676.5
$ws.Range("L40").Value = 4804.8203
$ws.Range("M40").Value = -2501.5
$ws.Range("N40").Value = -5154.8203

$ws.Range("H52").Value = 4651.077
$ws.Range("I52").Value = 2748
$ws.Range("J52").Value = 4997.091
$ws.Range("K52").Value = 8244
$ws.Range("L52").Value = 14991.273
$ws.Range("M52").Value = -8084
$ws.Range("N52").Value = -15311.273

$ws.Range("H53").Value = 209
$ws.Range("J53").Value = 201.6
$ws.Range("L53").Value = 201.6
$ws.Range("N53").Value = -1475.6

$ws.Range("H107").Value = 3030.7144
$ws.Range("I107").Value = 3749.3125
$ws.Range("K107").Value = 3749.3125
$ws.Range("M107").Value = -1829.3125

$ws.Range("H121").Value = 4949.8
$ws.Range("J121").Value = 4949.8
$ws.Range("L121").Value = 14849.4
$ws.Range("N121").Value = -18343.4

$ws.Range("H125").Value = 4467.467
$ws.Range("I125").Value = 4429.5
$ws.Range("K125").Value = 39865.5
$ws.Range("M125").Value = -37405.5

$ws.Range("H129").Value = 1744.7
$ws.Range("I129").Value = 1289.4
$ws.Range("J129").Value = 2200
$ws.Range("K129").Value = 3868.2
$ws.Range("L129").Value = 6600
$ws.Range("M129").Value = 1131.8
$ws.Range("N129").Value = -16600

$ws.Range("H131").Value = 2888
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

$ws.Range("H132").Value = 21018.4
$ws.Range("I132").Value = 26023
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 78069
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -75539
$ws.Range("N132").Value = -8060

$ws.Range("H137").Value = 20841016
$ws.Range("I137").Value = 38463364
$ws.Range("K137").Value = 115390092
$ws.Range("M137").Value = -115387542

$ws.Range("H141").Value = 8533.267
$ws.Range("I141").Value = 4222.1113
$ws.Range("K141").Value = 12666.3339
$ws.Range("M141").Value = -7486.333899999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 44500
$ws.Range("I76").Value = 45000
$ws.Range("J76").Value = 44000
$ws.Range("K76").Value = 45000
$ws.Range("L76").Value = 44000
$ws.Range("M76").Value = -44662
$ws.Range("N76").Value = -44676

$ws.Range("H79").Value = 44500
$ws.Range("I79").Value = 45000
$ws.Range("J79").Value = 44000
$ws.Range("K79").Value = 45000
$ws.Range("L79").Value = 44000
$ws.Range("M79").Value = -43830
$ws.Range("N79").Value = -46340

$ws.Range("H132").Value = 5078.6665
$ws.Range("J132").Value = 6646.2666
$ws.Range("L132").Value = 19938.7998
$ws.Range("N132").Value = -24998.7998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6580683
$ws.Range("I31").Value = 7520666
$ws.Range("J31").Value = 799
$ws.Range("K31").Value = 7520666
$ws.Range("L31").Value = 799
$ws.Range("M31").Value = -7520371
$ws.Range("N31").Value = -1389

$ws.Range("H34").Value = 6580683
$ws.Range("I34").Value = 7520666
$ws.Range("J34").Value = 799
$ws.Range("K34").Value = 7520666
$ws.Range("L34").Value = 799
$ws.Range("M34").Value = -7520464
$ws.Range("N34").Value = -1203

$ws.Range("H58").Value = 9400357
$ws.Range("I58").Value = 27782028
$ws.Range("J58").Value = 3273134.2
$ws.Range("K58").Value = 27782028
$ws.Range("L58").Value = 3273134.2
$ws.Range("M58").Value = -27781825
$ws.Range("N58").Value = -3273540.2

$ws.Range("H62").Value = 905
$ws.Range("I62").Value = 905
$ws.Range("K62").Value = 905
$ws.Range("M62").Value = -281

$ws.Range("H65").Value = 905
$ws.Range("I65").Value = 905
$ws.Range("K65").Value = 4525
$ws.Range("M65").Value = -1405

$ws.Range("H86").Value = 37791.867
$ws.Range("I86").Value = 78833.73
$ws.Range("K86").Value = 78833.73
$ws.Range("M86").Value = -77710.73

$ws.Range("H89").Value = 37791.867
$ws.Range("I89").Value = 78833.73
$ws.Range("K89").Value = 394168.65
$ws.Range("M89").Value = -388552.65

$ws.Range("H122").Value = 44299.668
$ws.Range("I122").Value = 2792
$ws.Range("K122").Value = 8376
$ws.Range("M122").Value = -5926

$ws.Range("H136").Value = 9400357
$ws.Range("I136").Value = 27782028
$ws.Range("J136").Value = 3273134.2
$ws.Range("K136").Value = 83346084
$ws.Range("L136").Value = 9819402.600000001
$ws.Range("M136").Value = -83343534
$ws.Range("N136").Value = -9824502.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 45454.547
$ws.Range("J37").Value = 45454.547
$ws.Range("L37").Value = 136363.641
$ws.Range("N37").Value = -136587.641

$ws.Range("H75").Value = 140
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

$ws.Range("H78").Value = 140
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

$ws.Range("H131").Value = 6598.8276
$ws.Range("I131").Value = 1941.5714
$ws.Range("J131").Value = 8080.6816
$ws.Range("K131").Value = 5824.7142
$ws.Range("L131").Value = 24242.0448
$ws.Range("M131").Value = -784.7142000000003
$ws.Range("N131").Value = -34322.0448

$ws.Range("H132").Value = 2147.6667
$ws.Range("I132").Value = 2147.6667
$ws.Range("K132").Value = 19329.0003
$ws.Range("M132").Value = -16799.0003

$ws.Range("H139").Value = 8031.5
$ws.Range("I139").Value = 4635.8335
$ws.Range("J139").Value = 13125
$ws.Range("K139").Value = 13907.5005
$ws.Range("L139").Value = 39375
$ws.Range("M139").Value = -8767.500499999998
$ws.Range("N139").Value = -49655

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 979.75
$ws.Range("I102").Value = 979.75
$ws.Range("K102").Value = 979.75
$ws.Range("M102").Value = 642.25

$ws.Range("H132").Value = 16410.867
$ws.Range("I132").Value = 13930.556
$ws.Range("K132").Value = 41791.66800000001
$ws.Range("M132").Value = -39261.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3155.7778
$ws.Range("I93").Value = 2252.6667
$ws.Range("K93").Value = 2252.6667
$ws.Range("M93").Value = -1004.6667

$ws.Range("H122").Value = 5022.6787
$ws.Range("I122").Value = 3113.2856
$ws.Range("J122").Value = 6932.0713
$ws.Range("K122").Value = 9339.856800000001
$ws.Range("L122").Value = 20796.2139
$ws.Range("M122").Value = -6889.856800000001
$ws.Range("N122").Value = -25696.2139

$ws.Range("H136").Value = 25003428
$ws.Range("I136").Value = 13892698
$ws.Range("K136").Value = 41678094
$ws.Range("M136").Value = -41675544

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2188
$ws.Range("I107").Value = 2199.2856
$ws.Range("K107").Value = 6597.8568
$ws.Range("M107").Value = -4677.8568

$ws.Range("H132").Value = 6412399
$ws.Range("I132").Value = 7248416.5
$ws.Range("J132").Value = 2932.3333
$ws.Range("K132").Value = 21745249.5
$ws.Range("L132").Value = 8796.999899999999
$ws.Range("M132").Value = -21742719.5
$ws.Range("N132").Value = -13856.9999

$ws.Range("H136").Value = 3569379
$ws.Range("I136").Value = 1614069.1
$ws.Range("K136").Value = 4842207.300000001
$ws.Range("M136").Value = -4839657.300000001
